$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so purely numeric-looking
# strings (e.g. "242.89") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '36.452.55'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.943.37'
$ws.Range("E3").Value = '  -2.06%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '242.89'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  -2.40%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '57.47'
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("E9").Value = '  -3.34%  '
$ws.Range("D10").Value = '0.0856'
$ws.Range("E10").Value = '  +4.28%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = '2.229.83'
$ws.Range("E12").Value = '  -2.02%  '
$ws.Range("E13").Value = '  -5.22%  '
$ws.Range("E14").Value = '  -9.68%  '
$ws.Range("E15").Value = '  -3.42%  '
$ws.Range("E16").Value = '  -4.83%  '
$ws.Range("D17").Value = '1.940.78'
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").Value = '36.362.50'
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").Value = '0.0₃0876'
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("D20").Value = '69.40'
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").Value = '228.71'
$ws.Range("E21").Value = '  -2.09%  '
$ws.Range("D22").Value = '5.01'
$ws.Range("E22").Value = '  -5.88%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  -7.11%  '
$ws.Range("D25").Value = '2.27'
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("D26").Value = '9.25'
$ws.Range("E26").Value = '  -8.23%  '
$ws.Range("D27").Value = '161.15'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("D28").Value = '0.130'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '19.30'
$ws.Range("E29").Value = '  -2.66%  '
$ws.Range("E30").Value = '  -1.90%  '
$ws.Range("D31").Value = '1.13'
$ws.Range("E31").Value = '  -6.06%  '
$ws.Range("E32").Value = '  -5.73%  '
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").Value = '4.24'
$ws.Range("E34").Value = '  -4.12%  '
$ws.Range("D35").Value = '6.16'
$ws.Range("E35").Value = '  -2.44%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("D39").Value = '3.08'
$ws.Range("E39").Value = '  +0.99%  '
$ws.Range("D40").Value = '0.0978'
$ws.Range("E40").Value = '  +2.04%  '
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("E42").Value = '  -6.08%  '
$ws.Range("E43").Value = '  -1.14%  '
$ws.Range("D44").Value = '15.77'
$ws.Range("E44").Value = '  -2.69%  '
$ws.Range("D45").Value = '1.344.66'
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("E46").Value = '  -6.30%  '
$ws.Range("D47").Value = '87.33'
$ws.Range("E47").Value = '  -5.65%  '
$ws.Range("E48").Value = '  -4.76%  '
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("E50").Value = '  -2.40%  '
$ws.Range("D51").Value = '2.120.32'

# Restore the original (unstyled) cell style now that the text values are set,
# so the D column cells keep matching the workbook's original "no explicit style" state.
$ws.Range("D2:D51").Style = "Normal"
